$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.360.85'
$ws.Range("E2").Value = '  -1.13%  '
$ws.Range("D3").Value = '2.435.98'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.32%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("D9").Value = '2.423.49'
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("E10").Value = '  +2.07%  '
$ws.Range("E11").Value = '  +1.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.14'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.340'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").Value = '2.891.52'
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").Value = '61.399.15'
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("D18").Value = '2.421.27'
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.39%  '
$ws.Range("E22").Value = '  -1.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.53%  '
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("E25").Value = '  -1.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '574.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.79%  '
$ws.Range("D29").Value = '2.569.05'
$ws.Range("E29").Value = '  +0.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").Value = '0.0₃0913'
$ws.Range("E31").Value = '  -3.03%  '
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("E33").Value = '  -4.88%  '
$ws.Range("E34").Value = '  -0.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.132'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.50%  '
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.62'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '152.02'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.52%  '
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("E40").Value = '  -3.17%  '
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.12'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.70'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.63%  '
$ws.Range("E45").Value = '  -5.03%  '
$ws.Range("E46").Value = '  -4.73%  '
$ws.Range("D47").Value = '0.0₆0288'
$ws.Range("E47").Value = '  +23.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '141.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.593'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0507'
$ws.Range("D51").Style = "Normal"
